# Trade #6 closed at 2026-02-17 19:44:16 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Sheet "Summary" ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1299.97   # Current Capital
$summary.Range("B4").Value = -0.03     # Total P&L $
$summary.Range("B5").Value = -0.1      # Total P&L %
$summary.Range("B6").Value = 6         # Total Trades
$summary.Range("B8").Value = 3         # Losing Trades
$summary.Range("B9").Value = 50        # Win Rate %

# --- Sheet "Strategy Status" (MarketMaking row) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.97      # Capital
$status.Range("D4").Value = 6          # Trades
$status.Range("E4").Value = -0.03      # P&L $
$status.Range("F4").Value = -0.03      # P&L %
$status.Range("G4").Value = 50         # Win Rate %

# --- Sheets "All Trades" and "MarketMaking" - append trade #6 row ---
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(7, 1).Value = 6                 # A7 Trade #

    # B7 is a date-looking string; force it to stay text instead of
    # being auto-converted to a date serial by Excel's input parser.
    $dateCell = $ws.Cells.Item(7, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.Style = "Normal"

    $ws.Cells.Item(7, 3).Value = "19:44:10"        # C7 Time
    $ws.Cells.Item(7, 4).Value = "MarketMaking"    # D7 Strategy
    $ws.Cells.Item(7, 5).Value = "UP"              # E7 Side
    $ws.Cells.Item(7, 6).Value = 0.22              # F7 Entry Price
    $ws.Cells.Item(7, 7).Value = 0.1               # G7 Exit Price
    $ws.Cells.Item(7, 8).Value = "CLOSED"          # H7 Status
    $ws.Cells.Item(7, 9).Value = -54.5455          # I7 P&L %
    $ws.Cells.Item(7, 10).Value = -0.12            # J7 P&L $
    $ws.Cells.Item(7, 11).Value = 99.97            # K7 Capital After
    $ws.Cells.Item(7, 12).Value = 0                # L7 Entry Slippage (bps)
    $ws.Cells.Item(7, 13).Value = 0                # M7 Exit Slippage (bps)
    $ws.Cells.Item(7, 14).Value = 0.6              # N7 Confidence
    $ws.Cells.Item(7, 15).Value = "Normal spread capture: 19600 bps"  # O7 Entry Reason
    $ws.Cells.Item(7, 16).Value = "early_exit"     # P7 Exit Reason
    $ws.Cells.Item(7, 17).Value = 0.14             # Q7 Duration (min)
}
